# Update "想去人数" (interest count) values in column F across sheets
# per the generated-output refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 143
$ws1.Range("F8").Value = 712
$ws1.Range("F18").Value = 412
$ws1.Range("F20").Value = 2765
$ws1.Range("F23").Value = 163
$ws1.Range("F26").Value = 955
$ws1.Range("F27").Value = 310
$ws1.Range("F28").Value = 241
$ws1.Range("F29").Value = 1290

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 149

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 143
$ws4.Range("F16").Value = 712
$ws4.Range("F21").Value = 149
$ws4.Range("F31").Value = 412
$ws4.Range("F33").Value = 2765
$ws4.Range("F36").Value = 163
$ws4.Range("F40").Value = 956
$ws4.Range("F41").Value = 310
$ws4.Range("F42").Value = 242
$ws4.Range("F43").Value = 1290
